$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows: rows 2 and 3 (columns A-E), matching the "V model" /
# secondary-write data dump. Column A holds date-serial values, B:E hold 1s.
$ws.Cells.Item(2, 1).Value = 35932
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = 1

$ws.Cells.Item(3, 1).Value = 3048
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 5).Value = 1

# --- Register the date/time number format + a plain border in the style
# table (mirrors the new numFmtId 22 + extra border style seen after the
# edit) via a scratch cell that is removed again, so the workbook ends up
# with the style definitions available without changing any visible cell.
$scratch = $ws.Cells.Item(50, 50)
$scratch.Value = 1
$scratch.NumberFormat = "m/d/yy h:mm"
$scratch.Borders.LineStyle = 1
$scratch.Borders.LineStyle = -4142
$scratch.EntireRow.Delete()

# --- Column widths: narrow "calendar strip" layout (mirrors the
# column width pattern seen in the target workbook).
$wide = 5.833333333333334
$narrow = 1.833333333333333

$ws.Columns.Item(1).ColumnWidth = $wide
$ws.Range("B1:I1").EntireColumn.ColumnWidth = $narrow
$ws.Columns.Item(10).ColumnWidth = $wide
$ws.Range("K1:N1").EntireColumn.ColumnWidth = $narrow

# --- Force a full recalculation on next load (mirrors calcPr/fullCalcOnLoad).
$wb.ForceFullCalculation = $true

Write-Host "Edit applied: wrote rows 2-3 (A:E), registered date style, resized columns."
